$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 38
$ws.Range("F4").Value = 1179
$ws.Range("F5").Value = 9264
$ws.Range("F6").Value = 148
$ws.Range("F7").Value = 254
$ws.Range("F8").Value = 7204
$ws.Range("F9").Value = 190
$ws.Range("F12").Value = 77
$ws.Range("F13").Value = 6486
$ws.Range("F15").Value = 443
$ws.Range("F16").Value = 428
$ws.Range("F17").Value = 623
$ws.Range("F19").Value = 0
$ws.Range("F20").Value = 153
$ws.Range("F22").Value = 162
$ws.Range("F23").Value = 10474
$ws.Range("F25").Value = 19
$ws.Range("F26").Value = 1995
$ws.Range("F27").Value = 0
$ws.Range("F29").Value = 2258
$ws.Range("F30").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("F32").Value = 193
$ws.Range("F34").Value = 2157
$ws.Range("F36").Value = 0
$ws.Range("F38").Value = 5379
$ws.Range("F39").Value = 443
$ws.Range("F40").Value = 1214
$ws.Range("F42").Value = 0
$ws.Range("F43").Value = 169
$ws.Range("F44").Value = 0
$ws.Range("F46").Value = 1008
$ws.Range("F47").Value = 1409
$ws.Range("F49").Value = 1104
# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 40
$ws.Range("F3").Value = 2
$ws.Range("F7").Value = 241
$ws.Range("F9").Value = 47
$ws.Range("F10").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = 9
$ws.Range("F14").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("F17").Value = 8
$ws.Range("F18").Value = 912
$ws.Range("F19").Value = 8
$ws.Range("F24").Value = 4
# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 0
# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 40
$ws.Range("F4").Value = 1179
$ws.Range("F5").Value = 9264
$ws.Range("F6").Value = 7205
$ws.Range("F7").Value = 190
$ws.Range("F8").Value = 73
$ws.Range("F9").Value = 6
$ws.Range("F12").Value = 5650
$ws.Range("F14").Value = 6486
$ws.Range("F15").Value = 6486
$ws.Range("F16").Value = 1107
$ws.Range("F17").Value = 443
$ws.Range("F18").Value = 428
$ws.Range("F26").Value = 0
$ws.Range("F27").Value = 10474
$ws.Range("F28").Value = 0
$ws.Range("F29").Value = 2337
$ws.Range("F30").Value = 0
$ws.Range("F34").Value = 193
$ws.Range("F35").Value = 26
$ws.Range("F36").Value = 2157
$ws.Range("F38").Value = 0
$ws.Range("F40").Value = 5379
$ws.Range("F41").Value = 443
$ws.Range("F42").Value = 0
$ws.Range("F43").Value = 721
$ws.Range("F44").Value = 129
$ws.Range("F45").Value = 0
$ws.Range("F46").Value = 1109
$ws.Range("F48").Value = 0
$ws.Range("F49").Value = 1409
$ws.Range("F50").Value = 69
$ws.Range("F51").Value = 1104
